$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift the existing extr1..extr8 rows (rows 8-15) down by two rows
# (to rows 10-17), going bottom-up so source data isn't clobbered before
# it is read, to make room for two new "line7"/"line8" rows.
for ($r = 15; $r -ge 8; $r--) {
    $src = $ws.Range("A" + $r + ":E" + $r)
    $dst = $ws.Range("A" + ($r + 2) + ":E" + ($r + 2))
    $src.Copy($dst) | Out-Null
}

# Renumber the shifted rows' index column (A) and flip in_service (E) for
# the first two of them (formerly extr1/extr2, now at rows 10 and 11)
# from false to true; all other shifted rows keep their original values.
$idx = 8
for ($r = 10; $r -le 17; $r++) {
    $ws.Cells.Item($r, 1).Value = $idx
    $idx++
}
$ws.Cells.Item(10, 5).Value = $true
$ws.Cells.Item(11, 5).Value = $true

# Write the new "line7" row (row 8), reusing row 7's formatting for column A.
$ws.Cells.Item(8, 1).Value = 6
$ws.Cells.Item(8, 2).Value = "line7"
$ws.Cells.Item(8, 3).Value = 14
$ws.Cells.Item(8, 4).Value = 11
$ws.Cells.Item(8, 5).Value = $true
$ws.Range("A7").Copy() | Out-Null
$ws.Cells.Item(8, 1).PasteSpecial(-4122) | Out-Null

# Write the new "line8" row (row 9), reusing row 7's formatting for column A.
$ws.Cells.Item(9, 1).Value = 7
$ws.Cells.Item(9, 2).Value = "line8"
$ws.Cells.Item(9, 3).Value = 16
$ws.Cells.Item(9, 4).Value = 9
$ws.Cells.Item(9, 5).Value = $false
$ws.Range("A7").Copy() | Out-Null
$ws.Cells.Item(9, 1).PasteSpecial(-4122) | Out-Null
